$wb = $excel.ActiveWorkbook

# "Chart" sheet holds the daily HTTPS/Non-HTTPS counts (A: Date, B: Non-HTTPS URLs, C: HTTPS URLs)
$chart = $wb.Worksheets.Item("Chart")

# Append two new days of data after the existing last row (row 36).
# Column A holds dates as plain text, matching the existing rows above.
# Force text storage via a "Text" number format so Excel does not
# auto-convert the "yyyy-MM-dd" strings into date serial numbers, then
# clear the formatting again so the cells keep the sheet's default style
# (matching the rest of the column).
$chart.Range("A37:A38").NumberFormat = "@"

$chart.Cells.Item(37, 1).Value = "2025-11-11"
$chart.Cells.Item(37, 2).Value = 0.0
$chart.Cells.Item(37, 3).Value = 54.0

$chart.Cells.Item(38, 1).Value = "2025-11-12"
$chart.Cells.Item(38, 2).Value = 0.0
$chart.Cells.Item(38, 3).Value = 46.0

$chart.Range("A37:A38").ClearFormats()

$wb.Save()
